$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update username/password values on rows 2 and 3
$ws.Range("A2").Value = "mngr515202"
$ws.Range("B2").Value = "usYrumy"
$ws.Range("A3").Value = "mngr515202"
$ws.Range("B3").Value = "usYrumy"

# Remove row 4 (duplicate of row 3) entirely
$ws.Rows.Item(4).Delete()
